# Move the last data row of Sheet1 (AFS / POSTED) out to a brand-new
# Sheet2, and update the remembered selections on both sheets to match
# the target workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1 loses its 5th row (AFS / POSTED) -> dimension becomes A1:B4.
$ws1.Rows(5).Delete() | Out-Null

# Remembered cursor position on Sheet1 moves to G18.
$ws1.Range("G18").Select() | Out-Null

# Add the new sheet right after Sheet1 -- it becomes "Sheet2" automatically.
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Re-create the row that was removed from Sheet1.
$ws2.Cells.Item(1, 1).Value = "AFS"
$ws2.Cells.Item(1, 2).Value = "POSTED"

# Remembered selection on Sheet2 covers the single data row.
$ws2.Range("A1:B1").Select() | Out-Null

# Keep Sheet1 as the active/visible tab (adding Sheet2 made it active).
$ws1.Activate() | Out-Null
